$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44475
$ws.Range("M2").Value = 240
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("S2").Value = 5750

$ws.Range("D4").Value = 44482
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10500
$ws.Range("S4").Value = 5250

$ws.Range("D5").Value = 44490
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 9500
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9750
$ws.Range("S5").Value = 4875

$ws.Range("D6").Value = 44455
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("S6").Value = 6250

$ws.Range("D7").Value = 44454
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12500
$ws.Range("S7").Value = 6250

$ws.Range("D8").Value = 44497
$ws.Range("M8").Value = 500
$ws.Range("N8").Value = 9000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 9500
$ws.Range("S8").Value = 4750

$ws.Range("D9").Value = 44517
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 5500
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 5750
$ws.Range("S9").Value = 2875

$ws.Range("D10").Value = 44489
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 9500
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 9750
$ws.Range("S10").Value = 4875
